$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the current row 1075 (Pehuenche / 44545 entry),
# pushing all rows from 1075 down to 1077+ (dimension grows from R1162 to R1164).
$ws.Rows.Item(1075).Insert()
$ws.Rows.Item(1075).Insert()

# New row 1075: Papa / Rodeo / 1a (guarda)
$ws.Cells.Item(1075, 1).Value = 10
$ws.Cells.Item(1075, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1075, 3).Value = "La Araucanía"
$ws.Cells.Item(1075, 4).Value = 45106
$ws.Cells.Item(1075, 5).Value = 9
$ws.Cells.Item(1075, 6).Value = 100114001
$ws.Cells.Item(1075, 7).Value = "Papa"
$ws.Cells.Item(1075, 8).Value = "Rodeo"
$ws.Cells.Item(1075, 9).Value = "1a (guarda)"
$ws.Cells.Item(1075, 10).Value = 280
$ws.Cells.Item(1075, 11).Value = 18000
$ws.Cells.Item(1075, 12).Value = 18000
$ws.Cells.Item(1075, 13).Value = 18000
$ws.Cells.Item(1075, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(1075, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(1075, 16).Value = 720
$ws.Cells.Item(1075, 17).Value = 25
$ws.Cells.Item(1075, 18).Value = "Hortaliza"

# New row 1076: Papa / Rosara / 1a (guarda)
$ws.Cells.Item(1076, 1).Value = 10
$ws.Cells.Item(1076, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(1076, 3).Value = "La Araucanía"
$ws.Cells.Item(1076, 4).Value = 45106
$ws.Cells.Item(1076, 5).Value = 9
$ws.Cells.Item(1076, 6).Value = 100114001
$ws.Cells.Item(1076, 7).Value = "Papa"
$ws.Cells.Item(1076, 8).Value = "Rosara"
$ws.Cells.Item(1076, 9).Value = "1a (guarda)"
$ws.Cells.Item(1076, 10).Value = 680
$ws.Cells.Item(1076, 11).Value = 15000
$ws.Cells.Item(1076, 12).Value = 15000
$ws.Cells.Item(1076, 13).Value = 15000
$ws.Cells.Item(1076, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(1076, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(1076, 16).Value = 600
$ws.Cells.Item(1076, 17).Value = 25
$ws.Cells.Item(1076, 18).Value = "Hortaliza"
